{"js": "// This change mirrors a regeneration of the document by a newer library\n// version (POI 3.17.0 -> 4.0.1). Only internal identifiers are affected:\n//   - the w:rsidR GUID stamped on the REF-field runs (fldChar/instrText/result)\n//   - the w:id used by the bookmarkStart/bookmarkEnd pair for \"bookmark1\"\n// No visible text or formatting changes.\n\n// Old values (for reference): rsidR \"5BEF7D3A56B04B6B848A298613E2C005\",\n// bookmark w:id \"113640858737380756001237403724904465710\".\nconst NEW_RSID = \"464D8D20643149FCAB03CCA7B1ABF91E\";\nconst NEW_BOOKMARK_ID = \"51102218167008688763084331361446262218\";\n// The source document uses a non-breaking space (U+00A0) before these colons.\nconst NBSP = \"\\u00A0\";\n\nconst NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapParagraphOoxml(paragraphXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document ' + NS + '><w:body>' + paragraphXml + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Paragraph 1: \"Test link before bookmark : \" + REF field to bookmark1\nconst beforeLinkOoxml =\n  '<w:p w:rsidP=\"009168BC\" w:rsidR=\"00E02A2B\" w:rsidRDefault=\"00E02A2B\">' +\n    '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Test link before bookmark' + NBSP + ': </w:t></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:instrText xml:space=\"preserve\"> REF bookmark1 \\\\h </w:instrText></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n  '</w:p>';\n\n// Paragraph 2: \"Test bookmark : \" + bookmarkStart/End bookmark1 + \"bookmarked content\"\nconst bookmarkParaOoxml =\n  '<w:p w:rsidP=\"00C31A62\" w:rsidR=\"00C31A62\" w:rsidRDefault=\"00C31A62\">' +\n    '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Test bookmark' + NBSP + ': </w:t></w:r>' +\n    '<w:bookmarkStart w:name=\"bookmark1\" w:id=\"' + NEW_BOOKMARK_ID + '\"/>' +\n    '<w:r><w:t>bookmarked content</w:t></w:r>' +\n    '<w:bookmarkEnd w:id=\"' + NEW_BOOKMARK_ID + '\"/>' +\n  '</w:p>';\n\n// Paragraph 4: \"Test link after bookmark : \" + REF field to bookmark1 + trailing space run\nconst afterLinkOoxml =\n  '<w:p w:rsidP=\"00E02A2B\" w:rsidR=\"00E02A2B\" w:rsidRDefault=\"00E02A2B\">' +\n    '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">Test link after bookmark' + NBSP + ': </w:t></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:instrText xml:space=\"preserve\"> REF bookmark1 \\\\h </w:instrText></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n    '<w:r w:rsidR=\"00D0546C\"><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '</w:p>';\n\n// NOTE: the document uses a non-breaking space (U+00A0) before the colons,\n// so we match on stable, ASCII-only substrings instead of the full text.\nfunction findParagraph(items, startsWith, mustContain) {\n  for (const p of items) {\n    if (p.text.indexOf(startsWith) === 0 && (!mustContain || p.text.indexOf(mustContain) !== -1)) {\n      return p;\n    }\n  }\n  return null;\n}\n\nconst beforePara = findParagraph(paragraphs.items, \"Test link before bookmark\");\nconst bookmarkPara = findParagraph(paragraphs.items, \"Test bookmark\", \"bookmarked content\");\nconst afterPara = findParagraph(paragraphs.items, \"Test link after bookmark\");\n\nif (!beforePara) throw new Error(\"Could not find 'Test link before bookmark' paragraph\");\nif (!bookmarkPara) throw new Error(\"Could not find 'Test bookmark : bookmarked content' paragraph\");\nif (!afterPara) throw new Error(\"Could not find 'Test link after bookmark' paragraph\");\n\nbeforePara.insertOoxml(wrapParagraphOoxml(beforeLinkOoxml), Word.InsertLocation.replace);\nbookmarkPara.insertOoxml(wrapParagraphOoxml(bookmarkParaOoxml), Word.InsertLocation.replace);\nafterPara.insertOoxml(wrapParagraphOoxml(afterLinkOoxml), Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# This change mirrors a regeneration of the document by a newer library\n# version (POI 3.17.0 -> 4.0.1). Only internal identifiers are affected:\n#   - the w:rsidR GUID stamped on the REF-field runs (fldChar/instrText/result)\n#   - the w:id used by the bookmarkStart/bookmarkEnd pair for \"bookmark1\"\n# No visible text or formatting changes.\n\n$d = $word.ActiveDocument\n\n$NEW_RSID = \"464D8D20643149FCAB03CCA7B1ABF91E\"\n$NEW_BOOKMARK_ID = \"51102218167008688763084331361446262218\"\n\n$xmlHeader = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Paragraph: \"Test link before bookmark : \" + REF field to bookmark1\n$beforeLinkXml = '<w:p w:rsidP=\"009168BC\" w:rsidR=\"00E02A2B\" w:rsidRDefault=\"00E02A2B\"><w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr><w:r><w:t xml:space=\"preserve\">Test link before bookmark&#160;: </w:t></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:instrText xml:space=\"preserve\"> REF bookmark1 \\h </w:instrText></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r></w:p>'\n\n# Paragraph: \"Test bookmark : \" + bookmarkStart/End bookmark1 + \"bookmarked content\"\n$bookmarkParaXml = '<w:p w:rsidP=\"00C31A62\" w:rsidR=\"00C31A62\" w:rsidRDefault=\"00C31A62\"><w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr><w:r><w:t xml:space=\"preserve\">Test bookmark&#160;: </w:t></w:r><w:bookmarkStart w:name=\"bookmark1\" w:id=\"' + $NEW_BOOKMARK_ID + '\"/><w:r><w:t>bookmarked content</w:t></w:r><w:bookmarkEnd w:id=\"' + $NEW_BOOKMARK_ID + '\"/></w:p>'\n\n# Paragraph: \"Test link after bookmark : \" + REF field to bookmark1 + trailing space run\n$afterLinkXml = '<w:p w:rsidP=\"00E02A2B\" w:rsidR=\"00E02A2B\" w:rsidRDefault=\"00E02A2B\"><w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\">Test link after bookmark&#160;: </w:t></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:instrText xml:space=\"preserve\"> REF bookmark1 \\h </w:instrText></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r><w:r w:rsidR=\"00D0546C\"><w:t xml:space=\"preserve\"> </w:t></w:r></w:p>'\n\n# Locate the three target paragraphs by their stable (ASCII) text prefixes -\n# the document text actually uses a non-breaking space (U+00A0) before the\n# colons, so we avoid relying on an exact full-text match.\n$beforeParaIndex = 0\n$bookmarkParaIndex = 0\n$afterParaIndex = 0\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($beforeParaIndex -eq 0 -and $t.StartsWith(\"Test link before bookmark\")) {\n        $beforeParaIndex = $i\n    }\n    if ($bookmarkParaIndex -eq 0 -and $t.StartsWith(\"Test bookmark\") -and $t.Contains(\"bookmarked content\") -and -not $t.Contains(\"duplicated\")) {\n        $bookmarkParaIndex = $i\n    }\n    if ($afterParaIndex -eq 0 -and $t.StartsWith(\"Test link after bookmark\")) {\n        $afterParaIndex = $i\n    }\n}\n\nif ($beforeParaIndex -eq 0) { throw \"Could not find 'Test link before bookmark' paragraph\" }\nif ($bookmarkParaIndex -eq 0) { throw \"Could not find 'Test bookmark : bookmarked content' paragraph\" }\nif ($afterParaIndex -eq 0) { throw \"Could not find 'Test link after bookmark' paragraph\" }\n\n[void]$d.Paragraphs.Item($beforeParaIndex).Range.InsertXML($xmlHeader + $beforeLinkXml + $xmlFooter)\n[void]$d.Paragraphs.Item($bookmarkParaIndex).Range.InsertXML($xmlHeader + $bookmarkParaXml + $xmlFooter)\n[void]$d.Paragraphs.Item($afterParaIndex).Range.InsertXML($xmlHeader + $afterLinkXml + $xmlFooter)\n"}
